# Trade #39 (MarketMaking) closed early at 2026-02-18 00:10:33; a new
# momentum trade (#68) opened at 2026-02-18 00:10:27. Updates the summary
# roll-ups, the per-strategy status row, the consolidated "All Trades"
# log (existing row + new row) and the strategy-specific "momentum" /
# "MarketMaking" sheets to match.

$wb = $excel.ActiveWorkbook

function Set-TextDate {
    param($cell, [string]$text)
    # Writing a yyyy-mm-dd-looking string through COM gets auto-parsed
    # into a date serial by Excel's type inference. Force Text format for
    # the write, then drop back to the default "Normal" style so the cell
    # ends up a plain text value with no lingering number format.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1) Summary sheet roll-up numbers
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1499.63
$wsSummary.Range("B4").Value = 0.73
$wsSummary.Range("B5").Value = 0.38
$wsSummary.Range("B6").Value = 38
$wsSummary.Range("B7").Value = 20
$wsSummary.Range("B9").Value = 52.63

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 6)
# ---------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C6").Value = 99.63
$wsStatus.Range("D6").Value = 9
$wsStatus.Range("E6").Value = -0.18
$wsStatus.Range("F6").Value = -0.37
$wsStatus.Range("G6").Value = 44.44

# ---------------------------------------------------------------------
# 3) All Trades sheet - close out trade #39 (row 40) + append trade #68
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

$wsAll.Range("G40").Value = 0.4
$wsAll.Range("H40").Value = "CLOSED"
$wsAll.Range("I40").Value = 17.6471
$wsAll.Range("J40").Value = 0.06
$wsAll.Range("K40").Value = 99.63
$wsAll.Range("L40").Value = "early_exit"
$wsAll.Range("M40").Value = 0.12

$wsAll.Range("A69").Value = 68
Set-TextDate $wsAll.Range("B69") "2026-02-18"
$wsAll.Range("C69").Value = "00:10:27"
$wsAll.Range("D69").Value = "momentum"
$wsAll.Range("E69").Value = "UP"
$wsAll.Range("F69").Value = 0.34
$wsAll.Range("G69").Value = ""
$wsAll.Range("H69").Value = "OPEN"
$wsAll.Range("I69").Value = 0
$wsAll.Range("J69").Value = 0
$wsAll.Range("K69").Value = 100
$wsAll.Range("L69").Value = ""
$wsAll.Range("M69").Value = 0
$wsAll.Range("N69").Value = 0
$wsAll.Range("O69").Value = 0
$wsAll.Range("P69").Value = 0.9
$wsAll.Range("Q69").Value = "Upward momentum: 21.687% over 10 samples"

# ---------------------------------------------------------------------
# 4) momentum sheet - append the same new trade #68 (row 6)
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

$wsMomentum.Range("A6").Value = 68
Set-TextDate $wsMomentum.Range("B6") "2026-02-18"
$wsMomentum.Range("C6").Value = "00:10:27"
$wsMomentum.Range("D6").Value = "momentum"
$wsMomentum.Range("E6").Value = "UP"
$wsMomentum.Range("F6").Value = 0.34
$wsMomentum.Range("G6").Value = ""
$wsMomentum.Range("H6").Value = "OPEN"
$wsMomentum.Range("I6").Value = 0
$wsMomentum.Range("J6").Value = 0
$wsMomentum.Range("K6").Value = 100
$wsMomentum.Range("L6").Value = 0
$wsMomentum.Range("M6").Value = 0
$wsMomentum.Range("N6").Value = 0.9
$wsMomentum.Range("O6").Value = "Upward momentum: 21.687% over 10 samples"
$wsMomentum.Range("P6").Value = ""
$wsMomentum.Range("Q6").Value = 0

# ---------------------------------------------------------------------
# 5) MarketMaking sheet - close out the same trade #39 (row 11)
# ---------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G11").Value = 0.4
$wsMM.Range("H11").Value = "CLOSED"
$wsMM.Range("I11").Value = 17.6471
$wsMM.Range("J11").Value = 0.06
$wsMM.Range("K11").Value = 99.63
$wsMM.Range("P11").Value = "early_exit"
$wsMM.Range("Q11").Value = 0.12

Write-Output "edit applied"
